# "Foreign Key Using email changing"
# Adds a "pencacah" (enumerator) column with the enumerator's email address
# as a foreign-key-style identifier, and moves the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column F: "pencacah"
$ws.Range("F1").Value = "pencacah"

# New data values in column F: enumerator email used as the foreign key
$ws.Range("F2").Value = "pcl01@bpssumsel.com"
$ws.Range("F3").Value = "pcl01@bpssumsel.com"

# The new email cells pick up an explicit (no-fill) formatting pass,
# same as the source edit's new cell style (xf index 4, applyFill).
$ws.Range("F2:F3").Interior.ColorIndex = -4142

# Move / collapse the current selection onto I11 (matches the saved view state)
$ws.Range("I11").Select()
